$wb = $excel.ActiveWorkbook

# Both the "展览" (Exhibition) and "全部类型" (All types) sheets contain the
# same event data and both need the "想去人数" (want-to-go count) column
# updated for three rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 282
    $ws.Range("F4").Value = 95
    $ws.Range("F5").Value = 849
}
